$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of Raw/Clean SSA data for 2020-07-18 after the existing last row (48)
$rDate = $ws.Range("A49")
$rDate.NumberFormat = "@"
$rDate.Value = "2020-07-18"
$rDate.ClearFormats()

$ws.Range("B49").Value = 338913
$ws.Range("C49").Value = 388636
$ws.Range("D49").Value = 87104
$ws.Range("E49").Value = 38888
$ws.Range("F49").Value = 28.66
